$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": drop three weekly entries ---
# Row 4 -> Order Week 2023-06-11 (45088.99999999999), qty 12
# Row 5 -> Order Week 2023-06-18 (45095.99999999999), qty 6
# Row 8 -> Order Week 2023-07-16 (45123.99999999999), qty 3
# Delete from the bottom up so earlier row numbers stay valid while deleting.
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows.Item(8).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# --- Sheet "Monthly Trend": adjust the two months whose weekly detail shrank ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(4, 2).Value = 6
$ws2.Cells.Item(5, 2).Value = 9
